# event_deck_15.xlsx -- "deploy: event logic & round reset"
# Adjusts a few event-card payout values, tweaks the header label, and
# resets the current selection back to the top of the card list (A11).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column header text tweak (A1, shared string).
$ws.Range("A1").Value = "イベン+A1:D22ト名"

# Event-logic value edits (payout columns H/I/M for a few cards).
$ws.Range("I5").Value = -3
$ws.Range("H11").Value = 2
$ws.Range("H12").Value = 1
$ws.Range("M12").Value = 1

# Round reset: move the selection back to A11 instead of the old A1:M21 block.
$ws.Range("A11").Select() | Out-Null
